# ValueSet-CDAContextControl.xlsx update for ST.r2b
#
# Changes applied (per commit "Update CDA Logical model for  ST.r2b"):
#  1. Rename the include sheet "Include from ContextControl" -> "Include #0".
#  2. On the Metadata sheet:
#       - bump the Version value.
#       - bump the Date value.
#       - insert a new "Jurisdiction" property row (with an empty value)
#         right after the "Contact" row, pushing Description/Purpose/
#         Copyright/Immutable down by one row.
#  3. The "Include #0" sheet's own property/value content is unchanged.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Update Version (row 3) and Date (row 8) values.
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" right after "Contact" (row 10), before
# "Description" (old row 11). Copy formatting from the row above so the new
# row matches the sheet's existing look (border/alignment/font), then set
# its own text.
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# Rename the second sheet.
$include = $wb.Worksheets.Item("Include from ContextControl")
$include.Name = "Include #0"
